$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new sheets, in final tab order, right after "Sheet1".
# ---------------------------------------------------------------------------
$sheet1 = $wb.Worksheets.Item(1)   # ChartData
$sheet2 = $wb.Worksheets.Item(2)   # Sheet1

$wsEntity       = $wb.Worksheets.Add($null, $sheet2)
$wsEntity.Name  = "EntityData"

$wsDateRange      = $wb.Worksheets.Add($null, $wsEntity)
$wsDateRange.Name = "DateRange"

$wsEngagement      = $wb.Worksheets.Add($null, $wsDateRange)
$wsEngagement.Name = "EngagementData"

$wsSearch      = $wb.Worksheets.Add($null, $wsEngagement)
$wsSearch.Name = "SearchData"

$wsSQL      = $wb.Worksheets.Add($null, $wsSearch)
$wsSQL.Name = "SQLData"

$wsDeleteSearch      = $wb.Worksheets.Add($null, $wsSQL)
$wsDeleteSearch.Name = "DeleteSearchData"

$wsShowEntry      = $wb.Worksheets.Add($null, $wsDeleteSearch)
$wsShowEntry.Name = "ShowEntryData"

# ---------------------------------------------------------------------------
# 2. Populate cell content.
#    The order of these writes matters: it controls the order new strings
#    are appended to the shared-string table, which must reproduce the
#    original authoring session.
# ---------------------------------------------------------------------------

# --- EntityData ---
$wsEntity.Range("A1").Value = "Entity "
$wsEntity.Range("A2").Value = "Account"
$wsEntity.Range("A3").Value = "Engagement"
$wsEntity.Range("A4").Value = "Project"
$wsEntity.Range("A5").Value = "Individual"

# --- DateRange (header filled in later, see below) ---
$wsDateRange.Range("A2").Value = "Today"
$wsDateRange.Range("A3").Value = "Yesterday"
$wsDateRange.Range("A4").Value = "Last 7 days"
$wsDateRange.Range("A5").Value = "Last 30 days"
$wsDateRange.Range("A6").Value = "This Month"
$wsDateRange.Range("A7").Value = "Last Month"
$wsDateRange.Range("A8").Value = "Custom"

# --- EngagementData header row ---
$wsEngagement.Range("A1").Value = "EngagementName"
$wsEngagement.Range("E1").NumberFormat = "@"
$wsEngagement.Range("E1").Value = "Years "
$wsEngagement.Range("C1").Value = "DeliveryMethod"
$wsEngagement.Range("D1").Value = "ContractType"
$wsEngagement.Range("B1").Value = "Account "

# --- SearchData ---
$wsSearch.Range("C2").Value = "AVA"
$wsSearch.Range("A1").Value = "Type"
$wsSearch.Range("B1").Value = "ID"
$wsSearch.Range("C1").Value = "Name"
$wsSearch.Range("A2").Value = "Engagement"
$wsSearch.Range("B2").Value = 1

# --- SQLData ---
$wsSQL.Range("A1").Value = "SQL search"
$wsSQL.Range("A2").Value = "SELECT * FROM accounts;"
$wsSQL.Range("A3").Value = "SELECT * FROM account"

# --- DeleteSearchData ---
$wsDeleteSearch.Range("A1").Value = "EntityID"
$wsDeleteSearch.Range("B1").Value = "EntityName"
$wsDeleteSearch.Range("A2").Value = 1
$wsDeleteSearch.Range("B2").Value = "Home TV"

# --- ShowEntryData header ---
$wsShowEntry.Range("A1").Value = "ShowEntry"

# --- DateRange header (typed after the other sheets were drafted) ---
$wsDateRange.Range("A1").Value = "DateRange"

# --- EngagementData data row ---
$wsEngagement.Range("C2").Value = "Utility"
$wsEngagement.Range("D2").Value = "Capacity Based"
$wsEngagement.Range("A2").Value = "Engagement1"
$wsEngagement.Range("B2").Value = "Affinion"
$wsEngagement.Range("E2").Value = 10
$wsEngagement.Range("E2").NumberFormat = "@"

# --- ShowEntryData data rows (stored as text) ---
$wsShowEntry.Range("A2").NumberFormat = "@"
$wsShowEntry.Range("A2").Value = "10"
$wsShowEntry.Range("A3").NumberFormat = "@"
$wsShowEntry.Range("A3").Value = "25"
$wsShowEntry.Range("A4").NumberFormat = "@"
$wsShowEntry.Range("A4").Value = "50"
$wsShowEntry.Range("A5").NumberFormat = "@"
$wsShowEntry.Range("A5").Value = "100"

# ---------------------------------------------------------------------------
# 3. Column widths.
# ---------------------------------------------------------------------------
$wsEntity.Columns.Item(1).ColumnWidth      = 15.022135416666666
$wsDateRange.Columns.Item(1).ColumnWidth   = 12.877604166666666
$wsEngagement.Columns.Item(1).ColumnWidth  = 17.166666666666668
$wsEngagement.Columns.Item(2).ColumnWidth  = 9.166666666666666
$wsEngagement.Columns.Item(3).ColumnWidth  = 14.736979166666666
$wsEngagement.Columns.Item(4).ColumnWidth  = 12.736979166666666
$wsSearch.Columns.Item(1).ColumnWidth      = 11.022135416666666
$wsSQL.Columns.Item(1).ColumnWidth         = 27.307291666666668
$wsDeleteSearch.Columns.Item(1).ColumnWidth = 10.166666666666666
$wsDeleteSearch.Columns.Item(2).ColumnWidth = 12.166666666666666

# ---------------------------------------------------------------------------
# 4. Page setup (portrait orientation) for the last two sheets.
# ---------------------------------------------------------------------------
$wsDeleteSearch.PageSetup.Orientation = 1
$wsShowEntry.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 5. Selections per sheet (also drives which sheet ends up "active").
# ---------------------------------------------------------------------------
$wsDateRange.Range("B4").Select() | Out-Null
$wsSearch.Range("E17").Select() | Out-Null
$wsSQL.Range("B6").Select() | Out-Null
$wsDeleteSearch.Range("H22").Select() | Out-Null
$wsShowEntry.Columns.Item(1).Select() | Out-Null

$sheet2.Range("B10").Select() | Out-Null

# EngagementData ends up active/selected last, matching the target workbook state.
$wsEngagement.Range("G14").Select() | Out-Null
